$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 2023-09-06 (45175) to 2023-09-08 (45177) for every data row (2-439).
$newDate = Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 439; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
